# Update the cryptos list with the latest scraped prices/volumes.
# Mirrors the automated "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.681.11"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "2.288.65"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'96.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "

$ws.Range("D6").Value = "'267.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.35%  "

$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.02%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("D10").Value = "'45.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.04%  "

$ws.Range("D11").Value = "'0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").Value = "'8.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("E13").Value = "  +0.12%  "

$ws.Range("D14").Value = "2.631.12"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").Value = "'15.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").Value = "'0.848"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "2.287.78"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "43.555.88"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("E19").Value = "  +2.41%  "

$ws.Range("E20").Value = "  -1.00%  "

$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("D22").Value = "'2.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.24%  "

$ws.Range("D23").Value = "'232.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.86%  "

$ws.Range("D24").Value = "'9.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.25%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").Value = "'2.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.68%  "

$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("E28").Value = "  +2.68%  "

$ws.Range("D29").Value = "'40.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.65%  "

$ws.Range("D31").Value = "'175.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("D32").Value = "'21.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.35%  "

$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("E34").Value = "  -3.95%  "

$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("E36").Value = "  -2.73%  "

$ws.Range("D37").Value = "'0.0353"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "

$ws.Range("E38").Value = "  -2.64%  "

$ws.Range("D39").Value = "'3.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("E40").Value = "  +2.34%  "

$ws.Range("D41").Value = "'2.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.88%  "

$ws.Range("D42").Value = "'12.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.33%  "

$ws.Range("D43").Value = "'65.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.86%  "

$ws.Range("E44").Value = "  +2.34%  "

$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("E46").Value = "  -4.90%  "

$ws.Range("D47").Value = "'0.101"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("D48").Value = "'97.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.78%  "

$ws.Range("E49").Value = "  -1.15%  "

$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "'0.187"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.72%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.511.96"
$ws.Range("E51").Value = "  -0.70%  "

